# Updates the EC (Estado de Cuenta) worksheet with a new worker-debt dataset
# and moves/reorders the signature footer, per the commit:
# "Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Remove the old worker rows that are no longer needed -------------
# Before: 19 data rows (16-34), a 4-row gap (35-38), then a 2-row footer (39-40).
# After:  8 data rows (16-23), the same 4-row gap (24-27), then the footer (28-29).
# Deleting rows 23-33 (11 rows) collapses the sheet so the old row 34 (which
# carries the special "closing" border style) becomes the new row 23, and the
# old footer rows 39-40 become the new rows 28-29, matching the target layout.
$ws.Range("B23:J33").EntireRow.Delete() | Out-Null

# --- 2. Update the summary figures at the top of the sheet ---------------
$ws.Range("E11").Value = 254708
$ws.Range("F13").Value = 7

# --- 3. Write the new worker debt detail table (rows 16-23) --------------
$data = @(
    @("CC", "1047448271", "JULIETH PAOLA JIMENEZ MALDONADO", "1902", 33125, 828116),
    @("CC", "1047448271", "JULIETH PAOLA JIMENEZ MALDONADO", "1903", 7729, 828116),
    @("CC", "1143367398", "JANER LATORRE SALCEDO", "2106", 20593, 908526),
    @("CC", "1017169516", "ANDRES FELIPE ARISTIZABAL GIRALDO", "2106", 36341, 908526),
    @("CC", "1017169516", "ANDRES FELIPE ARISTIZABAL GIRALDO", "2107", 36341, 908526),
    @("CC", "1001976349", "YAN CARLOS POLO CORPAS", "2408", 17333, 1423500),
    @("CC", "1001976349", "YAN CARLOS POLO CORPAS", "2409", 52000, 1423500),
    @("CC", "1007787040", "CARLOS ANDRES DIAZ MONTES", "2508", 51246, 1423500)
)

$row = 16
foreach ($rec in $data) {
    $ws.Range("B$row").Value = $rec[0]
    $ws.Range("C$row").Value = $rec[1]
    $ws.Range("D$row").Value = $rec[2]
    $ws.Range("E$row").Value = $rec[3]
    $ws.Range("F$row").Value = $rec[4]
    $ws.Range("G$row").Value = $rec[5]
    $row++
}

# --- 4. Reorder the signature footer --------------------------------------
# Previously the label row came first, then the signature line; now the
# signature line comes first, followed by the label row.
$ws.Range("B28").Value = "___________________________________"
$ws.Range("H28").Value = "___________________________________"

$ws.Range("B29").Value = "NOMBRE DEL REPRESENTANTE LEGAL"
$ws.Range("H29").Value = "FIRMA DEL REPRESENTANTE LEGAL"
